# =====================================================================
# npp-transect-info.xlsx — "add POC and methods edits"
#
# 1. Insert a new worksheet "ColumnHeadersPOC" (3rd tab, right after
#    "ColumnHeadersIntegrated") describing the new POC columns.
# 2. Add two new attribute rows (vol_filt_L / POC_ug_L) worth of
#    shared-string content (handled implicitly by writing the values).
# 3. Tidy up ColumnHeadersDiscrete: drop the heavy grey border that used
#    to sit under the last two rows of that sheet (incub_type block),
#    and remove the "thick bottom" row flag that went with it.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet in the correct position.
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("ColumnHeadersIntegrated")
$poc = $wb.Worksheets.Add($null, $afterSheet)
$poc.Name = "ColumnHeadersPOC"

# Column widths: A:G all ~37.11 characters wide.
$poc.Columns.Item("A:G").ColumnWidth = 37.109375

# ---------------------------------------------------------------------
# Header row.
# ---------------------------------------------------------------------
$poc.Range("A1").Value = "attributeName"
$poc.Range("B1").Value = "attributeDefinition"
$poc.Range("C1").Value = "class"
$poc.Range("D1").Value = "unit"
$poc.Range("E1").Value = "dateTimeFormatString"
$poc.Range("F1").Value = "missingValueCode"
$poc.Range("G1").Value = "missingValueCodeExplanation"

# ---------------------------------------------------------------------
# cruise
# ---------------------------------------------------------------------
$poc.Range("A2").Value = "cruise"
$poc.Range("B2").Value = "Identifier for research cruise generally including abbreviation for research vessel and voyage number"
$poc.Range("C2").Value = "character"
$poc.Range("A2:G2").WrapText = $true
$poc.Rows.Item(2).RowHeight = 43.2

# ---------------------------------------------------------------------
# date_time_utc
# ---------------------------------------------------------------------
$poc.Range("A3").Value = "date_time_utc"
$poc.Range("B3").Value = "Date and time in UTC when rosette bottle closed"
$poc.Range("C3").Value = "Date"
$poc.Range("E3").Value = "YYYY-MM-DD hh:mm:ss"
$poc.Range("A3:G3").WrapText = $true
$poc.Rows.Item(3).RowHeight = 28.8

# ---------------------------------------------------------------------
# latitude
# ---------------------------------------------------------------------
$poc.Range("A4").Value = "latitude"
$poc.Range("B4").Value = "Ship's latitude when rosette bottle closed"
$poc.Range("C4").Value = "numeric"
$poc.Range("D4").Value = "degree"
$poc.Range("B4:G4").WrapText = $true

# ---------------------------------------------------------------------
# longitude
# ---------------------------------------------------------------------
$poc.Range("A5").Value = "longitude"
$poc.Range("B5").Value = "Ship's longitude when rosette bottle closed"
$poc.Range("C5").Value = "numeric"
$poc.Range("D5").Value = "degree"
$poc.Range("B5").WrapText = $true

# ---------------------------------------------------------------------
# station
# ---------------------------------------------------------------------
$poc.Range("A6").Value = "station"
$poc.Range("B6").Value = "NES-LTER standard station from which sample was collected"
$poc.Range("C6").Value = "character"
$poc.Range("F6").Value = "NA"
$poc.Range("G6").Value = "Distance greater than 2 km from NES-LTER standard station"
$poc.Range("B6").Font.Color = 3355443
$poc.Range("B6").Font.Size = 12

# ---------------------------------------------------------------------
# cast
# ---------------------------------------------------------------------
$poc.Range("A7").Value = "cast"
$poc.Range("B7").Value = "CTD rosette cast number chronological per cruise"
$poc.Range("C7").Value = "numeric"
$poc.Range("D7").Value = "dimensionless"
$poc.Range("A7:G7").WrapText = $true
$poc.Rows.Item(7).RowHeight = 28.8

# ---------------------------------------------------------------------
# niskin
# ---------------------------------------------------------------------
$poc.Range("A8").Value = "niskin"
$poc.Range("B8").Value = "Rosette bottle position number"
$poc.Range("C8").Value = "numeric"
$poc.Range("D8").Value = "dimensionless"
$poc.Range("A8:G8").WrapText = $true

# ---------------------------------------------------------------------
# depth
# ---------------------------------------------------------------------
$poc.Range("A9").Value = "depth"
$poc.Range("B9").Value = "Depth of sample below sea surface http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"
$poc.Range("C9").Value = "numeric"
$poc.Range("D9").Value = "meter"
$poc.Range("B9").WrapText = $true
$poc.Rows.Item(9).RowHeight = 43.2

# ---------------------------------------------------------------------
# alternate_sample_category
# ---------------------------------------------------------------------
$poc.Range("A10").Value = "alternate_sample_category"
$poc.Range("B10").Value = "Identifer combining filter size with other information"
$poc.Range("C10").Value = "character"

# ---------------------------------------------------------------------
# filter_size
# ---------------------------------------------------------------------
$poc.Range("A11").Value = "filter_size"
$poc.Range("B11").Value = "Size fraction with filter and when applicable prefilter pore size in micrometers "
$poc.Range("C11").Value = "categorical"

# ---------------------------------------------------------------------
# vol_filt_L  (new attribute)
# ---------------------------------------------------------------------
$poc.Range("A12").Value = "vol_filt_L"
$poc.Range("B12").Value = "Volume of seawater filtered to obtain POC sample"
$poc.Range("C12").Value = "numeric"
$poc.Range("D12").Value = "liter"
$poc.Range("B12").WrapText = $true
$poc.Rows.Item(12).RowHeight = 28.8

# ---------------------------------------------------------------------
# POC_ug_L  (new attribute)
# ---------------------------------------------------------------------
$poc.Range("A13").Value = "POC_ug_L"
$poc.Range("B13").Value = "Particulate organic carbon http://vocab.nerc.ac.uk/collection/P09/current/POCP/"
$poc.Range("C13").Value = "numeric"
$poc.Range("D13").Value = "microgramsPerLiter"
$poc.Range("B13").WrapText = $true
$poc.Range("D13").Font.Color = 3355443
$poc.Range("D13").Font.Size = 12
$poc.Rows.Item(13).RowHeight = 43.2

$poc.Range("B13").Select()

# ---------------------------------------------------------------------
# 2. ColumnHeadersDiscrete — drop the grey border under the last block
#    (incub_type / ambient_temp rows) and the associated thick-bottom
#    row flag.
# ---------------------------------------------------------------------
$discrete = $wb.Worksheets.Item("ColumnHeadersDiscrete")

$lastRowVals = $discrete.Range("A17:D17").Value()
$discrete.Rows.Item(17).Delete()
$discrete.Rows.Item(16).Delete()

# Re-insert the two rows, now free of any stored thick-bottom flag.
$discrete.Rows.Item(16).Insert()
$discrete.Range("A16").Value = "iode_quality_flag"
$discrete.Range("B16").Value = "IODE Quality Flag primary level"
$discrete.Range("C16").Value = "categorical"
$discrete.Range("C16").WrapText = $true
$discrete.Range("F16").Value = "NA"
$discrete.Range("G16").Value = "Missing value"

$discrete.Rows.Item(17).Insert()
$discrete.Range("A17").Value = "incub_type"
$discrete.Range("B17").Value = "Identifier for incubation type (either ambient or experiemental)"
$discrete.Range("C17").Value = "categorical"
$discrete.Range("A17:D17").WrapText = $true

$discrete.Range("A6").Select()
